$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------
# The runs "{{ p" + "lot" + "sHeaderLabel }}" get merged into a single
# run "{{ plotsHeaderLabel }}". A simple Find/Replace across the split
# runs merges them into one run (using the first run's formatting),
# which matches the target exactly.
$d.Content.Find.Execute(
    "{{ plotsHeaderLabel }}", $false, $false, $false, $false, $false,
    $true, 1, $false, "{{ plotsHeaderLabel }}", 2) | Out-Null

# --- Change 2 -------------------------------------------------------
# The runs "{% if p" + "lot" + "sHeader%}{{ p" + "lots" + "Header }}{%
# else %}-{% endif%}" get merged into a single run
# "{% if plotsHeader%}{{ plotsHeader }}{% else %}-{% endif%}".
$d.Content.Find.Execute(
    "{% if plotsHeader%}{{ plotsHeader }}{% else %}-{% endif%}", $false,
    $false, $false, $false, $false, $true, 1, $false,
    "{% if plotsHeader%}{{ plotsHeader }}{% else %}-{% endif%}", 2) | Out-Null

# --- Change 3 -------------------------------------------------------
# The paragraph holding "{% if inputDateHeader %}..." is numbered at
# ilvl=5 (ListLevelNumber 6) and must move to ilvl=2 (ListLevelNumber
# 3). ListFormat.ListLevelNumber can't be assigned directly here, but
# ListOutdent() (the real Word list-level "Shift+Tab" operation) works
# and is idempotent-safe to call exactly as many times as needed.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "{% if inputDateHeader %}{{ inputDateHeader }}{% else %}-{% endif%}") {
        while ($p.Range.ListFormat.ListLevelNumber -gt 3) {
            $p.Range.ListFormat.ListOutdent()
        }
        break
    }
}

# --- Change 4 -------------------------------------------------------
# Split the single run
#   "{% if inputDateHeader %}{{ inputDateHeader }}{% else %}-{% endif%}"
# into three runs:
#   "{% if inputDateHeader %}{{ inputDateHeader }}"
#   "{% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}"
#   "{% else %}-{% endif %}"
$d.Content.Find.Execute(
    "{% if inputDateHeader %}{{ inputDateHeader }}{% else %}-{% endif%}",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "{% if inputDateHeader %}{{ inputDateHeader }}{% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}{% else %}-{% endif %}",
    2) | Out-Null
